$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 778.2421000000001
$ws.Range("I15").Value = 778.2421000000001
$ws.Range("K15").Value = 2334.7263
$ws.Range("M15").Value = -2165.7263

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2905.15
$ws.Range("I19").Value = 1676.8462
$ws.Range("K19").Value = 1676.8462
$ws.Range("M19").Value = -1501.8462

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3003.4688
$ws.Range("I76").Value = 2884.762
$ws.Range("J76").Value = 3230.0908
$ws.Range("K76").Value = 2884.762
$ws.Range("L76").Value = 3230.0908
$ws.Range("M76").Value = -2569.762
$ws.Range("N76").Value = -3860.0908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3003.4688
$ws.Range("I79").Value = 2884.762
$ws.Range("J79").Value = 3230.0908
$ws.Range("K79").Value = 2884.762
$ws.Range("L79").Value = 3230.0908
$ws.Range("M79").Value = -1792.762
$ws.Range("N79").Value = -5414.0908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1639.7097
$ws.Range("J112").Value = 1677.7
$ws.Range("L112").Value = 5033.1
$ws.Range("N112").Value = -7249.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3616.5715
$ws.Range("I132").Value = 2955.625
$ws.Range("J132").Value = 10666.667
$ws.Range("K132").Value = 8866.875
$ws.Range("L132").Value = 32000.001
$ws.Range("M132").Value = -6336.875
$ws.Range("N132").Value = -37060.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3772.2727
$ws.Range("I141").Value = 1184.2858
$ws.Range("J141").Value = 4980
$ws.Range("K141").Value = 3552.8574
$ws.Range("L141").Value = 14940
$ws.Range("M141").Value = 1627.1426
$ws.Range("N141").Value = -25300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 981.42426
$ws.Range("I2").Value = 885.875
$ws.Range("J2").Value = 1236.2222
$ws.Range("K2").Value = 885.875
$ws.Range("L2").Value = 1236.2222
$ws.Range("M2").Value = -772.875
$ws.Range("N2").Value = -1462.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 981.42426
$ws.Range("I116").Value = 885.875
$ws.Range("J116").Value = 1236.2222
$ws.Range("K116").Value = 885.875
$ws.Range("L116").Value = 1236.2222
$ws.Range("M116").Value = 1408.125
$ws.Range("N116").Value = -5824.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8841.861999999999
$ws.Range("I122").Value = 10517.941
$ws.Range("J122").Value = 6467.4165
$ws.Range("K122").Value = 31553.823
$ws.Range("L122").Value = 19402.2495
$ws.Range("M122").Value = -29103.823
$ws.Range("N122").Value = -24302.2495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 981.42426
$ws.Range("I3").Value = 885.875
$ws.Range("J3").Value = 1236.2222
$ws.Range("K3").Value = 885.875
$ws.Range("L3").Value = 1236.2222
$ws.Range("M3").Value = -771.875
$ws.Range("N3").Value = -1464.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 977.1429000000001
$ws.Range("I99").Value = 957.2727
$ws.Range("J99").Value = 1050
$ws.Range("K99").Value = 957.2727
$ws.Range("L99").Value = 1050
$ws.Range("M99").Value = 540.7273
$ws.Range("N99").Value = -4046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4209.9287
$ws.Range("I105").Value = 2100
$ws.Range("J105").Value = 4668.609
$ws.Range("K105").Value = 2100
$ws.Range("L105").Value = 4668.609
$ws.Range("M105").Value = -353
$ws.Range("N105").Value = -8162.609

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 7710.2563
$ws.Range("J118").Value = 7710.2563
$ws.Range("L118").Value = 7710.2563
$ws.Range("N118").Value = -11024.2563

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1640.625
$ws.Range("I99").Value = 1335.1666
$ws.Range("J99").Value = 2557
$ws.Range("K99").Value = 1335.1666
$ws.Range("L99").Value = 2557
$ws.Range("M99").Value = 162.8334
$ws.Range("N99").Value = -5553

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1640.625
$ws.Range("I126").Value = 1335.1666
$ws.Range("J126").Value = 2557
$ws.Range("K126").Value = 4005.4998
$ws.Range("L126").Value = 7671
$ws.Range("M126").Value = -1535.4998
$ws.Range("N126").Value = -12611

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 858.8333
$ws.Range("I5").Value = 300.70587
$ws.Range("J5").Value = 2214.2856
$ws.Range("K5").Value = 902.11761
$ws.Range("L5").Value = 6642.8568
$ws.Range("M5").Value = -790.11761
$ws.Range("N5").Value = -6866.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1641.8788
$ws.Range("J131").Value = 1254.28
$ws.Range("L131").Value = 3762.84
$ws.Range("N131").Value = -13842.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1160.2667
$ws.Range("I132").Value = 934
$ws.Range("J132").Value = 1311.1111
$ws.Range("K132").Value = 8406
$ws.Range("L132").Value = 11799.9999
$ws.Range("M132").Value = -5876
$ws.Range("N132").Value = -16859.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 858.8333
$ws.Range("I135").Value = 300.70587
$ws.Range("J135").Value = 2214.2856
$ws.Range("K135").Value = 2706.35283
$ws.Range("L135").Value = 19928.5704
$ws.Range("M135").Value = -171.3528299999998
$ws.Range("N135").Value = -24998.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11584
$ws.Range("I102").Value = 13976.5
$ws.Range("J102").Value = 2014
$ws.Range("K102").Value = 13976.5
$ws.Range("L102").Value = 2014
$ws.Range("M102").Value = -12354.5
$ws.Range("N102").Value = -5258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 14521.739
$ws.Range("J118").Value = 14521.739
$ws.Range("L118").Value = 14521.739
$ws.Range("N118").Value = -17835.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2780452.2
$ws.Range("I122").Value = 4446430
$ws.Range("J122").Value = 3822.2222
$ws.Range("K122").Value = 13339290
$ws.Range("L122").Value = 11466.6666
$ws.Range("M122").Value = -13336840
$ws.Range("N122").Value = -16366.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11948.923
$ws.Range("I132").Value = 19289.143
$ws.Range("J132").Value = 3385.3333
$ws.Range("K132").Value = 57867.429
$ws.Range("L132").Value = 10155.9999
$ws.Range("M132").Value = -55337.429
$ws.Range("N132").Value = -15215.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22739134
$ws.Range("I132").Value = 9574
$ws.Range("J132").Value = 38474984
$ws.Range("K132").Value = 28722
$ws.Range("L132").Value = 115424952
$ws.Range("M132").Value = -26192
$ws.Range("N132").Value = -115430012

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 29419384
$ws.Range("I136").Value = 35716224
$ws.Range("J136").Value = 34135
$ws.Range("K136").Value = 107148672
$ws.Range("L136").Value = 102405
$ws.Range("M136").Value = -107146122
$ws.Range("N136").Value = -107505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1539.091
$ws.Range("I113").Value = 204.28572
$ws.Range("J113").Value = 3875
$ws.Range("K113").Value = 612.85716
$ws.Range("L113").Value = 11625
$ws.Range("M113").Value = 1557.14284
$ws.Range("N113").Value = -15965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2549
$ws.Range("I122").Value = 2312.375
$ws.Range("K122").Value = 6937.125
$ws.Range("M122").Value = -4487.125
